# Implementação da restrição 4
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planilha1")

# Update the constraint coefficients (restrição 4): K11 from -1 to -2, M11 from 6 to 0
$ws.Range("K11").Value = -2
$ws.Range("M11").Value = 0

# Update the selected cell on the sheet
$ws.Activate()
$ws.Range("K12").Select()
